$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark that sits after "her large.  The file
#    features" (it gets re-added later, at the end of the document edit).
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()
Write-Output "removed old _GoBack bookmark"

# ---------------------------------------------------------------------------
# 2) Update the false-positive / false-negative rate table values.
# ---------------------------------------------------------------------------
$tableEdits = @{
    "0.0337" = "0.0556"
    "0.0738" = "0.1873"
    "0.0315" = "0.0521"
    "0.0772" = "0.1956"
    "0.0326" = "0.0539"
    "0.0554" = "0.1405"
    "0.0348" = "0.0573"
    "0.0522" = "0.1326"
    "0.0641" = "0.1630"
    "0.0335" = "0.0552"
    "0.0645" = "0.1638"
}

foreach ($old in $tableEdits.Keys) {
    $new = $tableEdits[$old]
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    Write-Output ("table cell " + $old + " -> " + $new)
}

# ---------------------------------------------------------------------------
# 3) Rewrite the "Although the false positives ... difference of ~6.5"
#    sentence, splitting it into several runs and updating the percentages.
# ---------------------------------------------------------------------------

# Locate each fixed sub-string in turn, always searching forward from the end
# of the previous piece, so the absolute character offsets are derived from
# the engine itself rather than hand computed.

$rng = $d.Content
$rng.Find.Execute("  Although the false positives are above 0%, they are relatively close")
$p0 = $rng.Start
$p1 = $rng.End
Write-Output ("run1 [" + $p0 + "," + $p1 + "] = [" + $rng.Text + "]")

$rng = $d.Range($p1, $d.Content.End)
$rng.Find.Execute(" ~3.3")
$p2old = $rng.End
Write-Output ("run2-old [" + $rng.Start + "," + $p2old + "] = [" + $rng.Text + "]")

# --- Run 2: " ~3.3" -> " at ~5.5"
$r2 = $d.Range($p1, $p2old)
$r2.Text = " at ~5.5"
$p2new = $r2.End
Write-Output ("run2 text replaced, new end=" + $p2new)

# mark run2 so it does not get coalesced back into run1/run3
$r2mark = $d.Range($p1, $p2new)
$r2mark.Bold = 1
Write-Output "run2 bold set"
$r2mark.Bold = 0
Write-Output "run2 bold cleared"

# --- Run 3: "%.  The real competitiv" (unchanged)
$rng = $d.Range($p2new, $d.Content.End)
$rng.Find.Execute("%.  The real competitiv")
$p3 = $rng.End
Write-Output ("run3 [" + $rng.Start + "," + $p3 + "] = [" + $rng.Text + "]")

# --- Run 4: "e" (unchanged, split off on its own)
$p4 = $p3 + 1
$r4 = $d.Range($p3, $p4)
Write-Output ("run4 text=[" + $r4.Text + "]")
$r4.Bold = 1
Write-Output "run4 bold set"
$r4.Bold = 0
Write-Output "run4 bold cleared"

# --- Run 5: " edge is the difference of ~6.5" -> " edge is the difference of ~16.4"
$rng = $d.Range($p4, $d.Content.End)
$rng.Find.Execute(" edge is the difference of ~6.5")
$p5old = $rng.End
Write-Output ("run5-old [" + $rng.Start + "," + $p5old + "] = [" + $rng.Text + "]")

$r5 = $d.Range($p4, $p5old)
$r5.Text = " edge is the difference of ~16.4"
$p5new = $r5.End
Write-Output ("run5 text replaced, new end=" + $p5new)

$r5mark = $d.Range($p4, $p5new)
Write-Output ("r5mark text=[" + $r5mark.Text + "]")
$r5mark.Bold = 1
Write-Output "run5 bold set"
$r5mark.Bold = 0
Write-Output "run5 bold cleared"

# ---------------------------------------------------------------------------
# 4) Re-insert the _GoBack bookmark right after the "~16.4" text (this is
#    where Word leaves it after the last edit made to the document).
# ---------------------------------------------------------------------------
$bmRange = $d.Range($p5new, $p5new)
$d.Bookmarks.Add("_GoBack", $bmRange)
Write-Output "added new _GoBack bookmark"
